$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.650772
$ws.Range("H2").Value = 163.952316
$ws.Range("I2").Value = 0.3089913429651142
$ws.Range("J2").Value = 0.3089913429651142
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 344.6193520131119
$ws.Range("R2").Value = 3101.574168118008
$ws.Range("S2").Value = 0.00419038767940246
$ws.Range("T2").Value = 0.004190387679402458
$ws.Range("G3").Value = 54.650772
$ws.Range("H3").Value = 163.952316
$ws.Range("I3").Value = 0.3089913429651142
$ws.Range("J3").Value = 0.3089913429651142
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 9975.497353975428
$ws.Range("R3").Value = 89779.47618577885
$ws.Range("S3").Value = 0.1212967320721444
$ws.Range("T3").Value = 0.1212967320721444
$ws.Range("G4").Value = 54.650772
$ws.Range("H4").Value = 163.952316
$ws.Range("I4").Value = 0.3089913429651142
$ws.Range("J4").Value = 0.3089913429651142
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 6962.311227465396
$ws.Range("R4").Value = 62660.80104718856
$ws.Range("S4").Value = 0.08465799444317444
$ws.Range("T4").Value = 0.08465799444317443
$ws.Range("G5").Value = 54.650772
$ws.Range("H5").Value = 163.952316
$ws.Range("I5").Value = 0.3089913429651142
$ws.Range("J5").Value = 0.3089913429651142
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 1061.536668195696
$ws.Range("R5").Value = 9553.830013761264
$ws.Range("S5").Value = 0.01290772021262444
$ws.Range("T5").Value = 0.01290772021262443
$ws.Range("G6").Value = 54.650772
$ws.Range("H6").Value = 163.952316
$ws.Range("I6").Value = 0.3089913429651142
$ws.Range("J6").Value = 0.3089913429651142
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 7067.621279464681
$ws.Range("R6").Value = 63608.59151518212
$ws.Range("S6").Value = 0.08593850855776847
$ws.Range("T6").Value = 0.08593850855776844
$ws.Range("I7").Value = 0.534813606173264
$ws.Range("J7").Value = 0.5348136061732639
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 596.4798775221174
$ws.Range("R7").Value = 5368.318897699057
$ws.Range("S7").Value = 0.007252877457923691
$ws.Range("T7").Value = 0.007252877457923688
$ws.Range("I8").Value = 0.534813606173264
$ws.Range("J8").Value = 0.5348136061732639
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.2099448550047561
$ws.Range("T8").Value = 0.2099448550047561
$ws.Range("I9").Value = 0.534813606173264
$ws.Range("J9").Value = 0.5348136061732639
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 12050.62491113794
$ws.Range("R9").Value = 108455.6242002415
$ws.Range("S9").Value = 0.146529177371361
$ws.Range("T9").Value = 0.146529177371361
$ws.Range("I10").Value = 0.534813606173264
$ws.Range("J10").Value = 0.5348136061732639
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 1837.346794751428
$ws.Range("R10").Value = 16536.12115276285
$ws.Range("S10").Value = 0.02234115793712899
$ws.Range("T10").Value = 0.02234115793712898
$ws.Range("I11").Value = 0.534813606173264
$ws.Range("J11").Value = 0.5348136061732639
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 12232.89943098554
$ws.Range("R11").Value = 110096.0948788699
$ws.Range("S11").Value = 0.1487455384020942
$ws.Range("T11").Value = 0.1487455384020942
$ws.Range("G12").Value = 5.807188
$ws.Range("H12").Value = 17.421564
$ws.Range("I12").Value = 0.0328334029567029
$ws.Range("J12").Value = 0.03283340295670289
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 36.619233221048
$ws.Range("R12").Value = 329.573098989432
$ws.Range("S12").Value = 0.0004452703622772943
$ws.Range("T12").Value = 0.0004452703622772942
$ws.Range("G13").Value = 5.807188
$ws.Range("H13").Value = 17.421564
$ws.Range("I13").Value = 0.0328334029567029
$ws.Range("J13").Value = 0.03283340295670289
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 1059.995795265945
$ws.Range("R13").Value = 9539.962157393509
$ws.Range("S13").Value = 0.01288898401890045
$ws.Range("T13").Value = 0.01288898401890045
$ws.Range("G14").Value = 5.807188
$ws.Range("H14").Value = 17.421564
$ws.Range("I14").Value = 0.0328334029567029
$ws.Range("J14").Value = 0.03283340295670289
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 739.8148046728841
$ws.Range("R14").Value = 6658.333242055956
$ws.Range("S14").Value = 0.008995753791629318
$ws.Range("T14").Value = 0.008995753791629316
$ws.Range("G15").Value = 5.807188
$ws.Range("H15").Value = 17.421564
$ws.Range("I15").Value = 0.0328334029567029
$ws.Range("J15").Value = 0.03283340295670289
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 112.7988274549173
$ws.Range("R15").Value = 1015.189447094256
$ws.Range("S15").Value = 0.001371573633508966
$ws.Range("T15").Value = 0.001371573633508966
$ws.Range("G16").Value = 5.807188
$ws.Range("H16").Value = 17.421564
$ws.Range("I16").Value = 0.0328334029567029
$ws.Range("J16").Value = 0.03283340295670289
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 751.0050449543868
$ws.Range("R16").Value = 6759.04540458948
$ws.Range("S16").Value = 0.009131821150386868
$ws.Range("T16").Value = 0.009131821150386864
$ws.Range("G17").Value = 15.46624066666667
$ws.Range("H17").Value = 46.398722
$ws.Range("I17").Value = 0.08744495821971184
$ws.Range("J17").Value = 0.08744495821971182
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 97.52773184293733
$ws.Range("R17").Value = 877.749586586436
$ws.Range("S17").Value = 0.001185885248542752
$ws.Range("T17").Value = 0.001185885248542752
$ws.Range("G18").Value = 15.46624066666667
$ws.Range("H18").Value = 46.398722
$ws.Range("I18").Value = 0.08744495821971184
$ws.Range("J18").Value = 0.08744495821971182
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 2823.078928258882
$ws.Range("R18").Value = 25407.71035432993
$ws.Range("S18").Value = 0.03432713540273451
$ws.Range("T18").Value = 0.03432713540273449
$ws.Range("G19").Value = 15.46624066666667
$ws.Range("H19").Value = 46.398722
$ws.Range("I19").Value = 0.08744495821971184
$ws.Range("J19").Value = 0.08744495821971182
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 1970.343274203249
$ws.Range("R19").Value = 17733.08946782924
$ws.Range("S19").Value = 0.02395832425597694
$ws.Range("T19").Value = 0.02395832425597694
$ws.Range("G20").Value = 15.46624066666667
$ws.Range("H20").Value = 46.398722
$ws.Range("I20").Value = 0.08744495821971184
$ws.Range("J20").Value = 0.08744495821971182
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 300.4162793309875
$ws.Range("R20").Value = 2703.746513978888
$ws.Range("S20").Value = 0.003652901870561818
$ws.Range("T20").Value = 0.003652901870561817
$ws.Range("G21").Value = 15.46624066666667
$ws.Range("H21").Value = 46.398722
$ws.Range("I21").Value = 0.08744495821971184
$ws.Range("J21").Value = 0.08744495821971182
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 2000.146158027838
$ws.Range("R21").Value = 18001.31542225054
$ws.Range("S21").Value = 0.02432071144189583
$ws.Range("T21").Value = 0.02432071144189582
$ws.Range("G22").Value = 6.352523666666666
$ws.Range("H22").Value = 19.057571
$ws.Range("I22").Value = 0.0359166896852071
$ws.Range("J22").Value = 0.03591668968520709
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 40.05803595335533
$ws.Range("R22").Value = 360.522323580198
$ws.Range("S22").Value = 0.0004870843710297915
$ws.Range("T22").Value = 0.0004870843710297914
$ws.Range("G23").Value = 6.352523666666666
$ws.Range("H23").Value = 19.057571
$ws.Range("I23").Value = 0.0359166896852071
$ws.Range("J23").Value = 0.03591668968520709
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 1159.536831938982
$ws.Range("R23").Value = 10435.83148745084
$ws.Range("S23").Value = 0.01409934998132549
$ws.Range("T23").Value = 0.01409934998132548
$ws.Range("G24").Value = 6.352523666666666
$ws.Range("H24").Value = 19.057571
$ws.Range("I24").Value = 0.0359166896852071
$ws.Range("J24").Value = 0.03591668968520709
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 809.2886015804677
$ws.Range("R24").Value = 7283.597414224209
$ws.Range("S24").Value = 0.009840518140764797
$ws.Range("T24").Value = 0.009840518140764795
$ws.Range("G25").Value = 6.352523666666666
$ws.Range("H25").Value = 19.057571
$ws.Range("I25").Value = 0.0359166896852071
$ws.Range("J25").Value = 0.03591668968520709
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 123.3914281713649
$ws.Range("R25").Value = 1110.522853542284
$ws.Range("S25").Value = 0.001500374013626165
$ws.Range("T25").Value = 0.001500374013626165
$ws.Range("G26").Value = 6.352523666666666
$ws.Range("H26").Value = 19.057571
$ws.Range("I26").Value = 0.0359166896852071
$ws.Range("J26").Value = 0.03591668968520709
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 821.5296838777746
$ws.Range("R26").Value = 7067.621279464681
$ws.Range("S26").Value = 0.009989363178460867
$ws.Range("T26").Value = 0.009989363178460864
